# Applies two edits to the document:
#  1. Merge the "F" / "fffffff" / "  " runs in the first paragraph into a
#     single run reading "ffffffff".
#  2. Change "Version 4" to "Version 5" in the third paragraph.

$d = $word.ActiveDocument

# --- Edit 1: collapse "Ffffffff  " (with trailing spaces) to "ffffffff" ---
$d.Content.Find.Execute(
    "Ffffffff  ",  # FindText
    $true,         # MatchCase
    $false,        # MatchWholeWord
    $false,        # MatchWildcards
    $false,        # MatchSoundsLike
    $false,        # MatchAllWordForms
    $true,         # Forward
    1,             # Wrap (wdFindContinue)
    $false,        # Format
    "ffffffff",    # ReplaceWith
    2              # Replace (wdReplaceAll)
)

# --- Edit 2: "Version 4" -> "Version 5" ---
$d.Content.Find.Execute(
    "Version 4",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Version 5",
    2
)
